$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 8
$ws.Cells.Item($row, 1).Value = 42604.890439814815
$ws.Cells.Item($row, 2).Value = "Named"
$ws.Cells.Item($row, 3).Value = 5765
$ws.Cells.Item($row, 4).Value = 2620
$ws.Cells.Item($row, 5).Value = 158
$ws.Cells.Item($row, 6).Value = 18
$ws.Cells.Item($row, 7).Value = 21
$ws.Cells.Item($row, 8).Value = 46
$ws.Cells.Item($row, 9).Value = 53
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 3
$ws.Cells.Item($row, 12).Value = 0
$ws.Cells.Item($row, 13).Value = 99
